# Updated symbol list on Sun Dec 25 09:52:20 UTC 2022 with GitHub Actions
#
# Applies the scraped-price refresh to Sheet1: most rows only get a new
# "Price" (column D) reading, a handful of rows (41-43) were re-ordered by
# the scraper (different coin landed on each row) so Coin/Link/Price/Volume
# all change together, and row 44 picked up a "Bestin24h" suffix on its
# Volume(1h) label alongside its new price.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text so that numeric-looking strings
# (e.g. "245.00", "0.006212") are preserved exactly instead of being
# auto-coerced into a floating point Number by Excel's type inference.
# The NumberFormat flip forces text entry; resetting the style back to
# "Normal" afterwards drops the now-unneeded text format override so the
# cell's style stays untouched, matching cells that were never reformatted.
function Set-TextValue($Cell, $Text) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

function Set-PlainValue($Cell, $Text) {
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Simple Price (column D) refreshes - coin identity/link/volume unchanged
# ---------------------------------------------------------------------
$priceUpdates = @{
    2  = "245.00"
    3  = "23.14"
    4  = "5.407"
    5  = "0.05981"
    6  = "3.389"
    7  = "0.8086"
    8  = "0.9274"
    9  = "0.1428"
    10 = "0.07431"
    11 = "0.03370"
    12 = "0.03036"
    13 = "0.09358"
    14 = "3.940"
    15 = "0.001593"
    16 = "0.04804"
    18 = "0.005731"
    19 = "0.004160"
    20 = "0.0009851"
    23 = "6.462"
    24 = "2.187"
    26 = "0.1313"
    40 = "0.03940"
    45 = "0.00005197"
    49 = "0.002237"
}

foreach ($row in $priceUpdates.Keys) {
    $addr = "D" + $row
    Set-TextValue $ws.Range($addr) $priceUpdates[$row]
}

# ---------------------------------------------------------------------
# Rows 41-43: the scraper's ranking reshuffled - each row now shows a
# different coin (with its own link/price/volume label).
# ---------------------------------------------------------------------

# Row 41: KickToken -> BKEXToken
Set-PlainValue $ws.Range("B41") "BKEXToken"
Set-PlainValue $ws.Range("C41") "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D41") "0.1076"
Set-PlainValue $ws.Range("E41") "40BKEXTokenBKK"

# Row 42: BKEXToken -> CEJI
Set-PlainValue $ws.Range("B42") "CEJI"
Set-PlainValue $ws.Range("C42") "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D42") "0.002711"
Set-PlainValue $ws.Range("E42") "41CEJICEJI"

# Row 43: CEJI -> KickToken
Set-PlainValue $ws.Range("B43") "KickToken"
Set-PlainValue $ws.Range("C43") "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D43") "0.006202"
Set-PlainValue $ws.Range("E43") "42KickTokenKICK"

# Row 44: LocalTraders keeps its spot, price ticks up and the volume
# label gains a "Bestin24h" marker.
Set-TextValue $ws.Range("D44") "0.007537"
Set-PlainValue $ws.Range("E44") "43LocalTradersLCTBestin24h"
